$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 0.59779730578045021
$ws.Range("BG1").Value = 0.63129596928858334
$ws.Range("D2").Value = 0.80337072208117632
$ws.Range("N2").Value = 0.6834968814075314
$ws.Range("R2").Value = 0.94746298503938653
$ws.Range("BP2").Value = 0.59631321955170424
$ws.Range("A3").Value = 0.78405351687826141
$ws.Range("BK3").Value = 0.61617900570898021
$ws.Range("BP3").Value = 0.99126420866485132
$ws.Range("M4").Value = 0.79230323661503865
$ws.Range("G6").Value = 0.95984167792591757
$ws.Range("AM6").Value = 0.78326645442414777
$ws.Range("BC6").Value = 0.69778017866738651
$ws.Range("L7").Value = 0.75094381675418398
$ws.Range("AF8").Value = 0.91838705368758999
$ws.Range("J9").Value = 0.99118797958950844
$ws.Range("AA9").Value = 0.97073205289174425
$ws.Range("AO9").Value = 0.61735360173977694
$ws.Range("AF10").Value = 0.78885428901921895
$ws.Range("AZ11").Value = 0.76240326314268592
$ws.Range("BP12").Value = 0.92408083034189026
$ws.Range("Q13").Value = 0.9470720371440704
$ws.Range("AU13").Value = 0.73794405119331186
$ws.Range("BE13").Value = 0.94885327931131902
$ws.Range("BK14").Value = 0.62469123496376155
$ws.Range("C15").Value = 0.80909279966027858
$ws.Range("E15").Value = 0.96728127051647372
$ws.Range("L15").Value = 0.89562882902271057
$ws.Range("AS16").Value = 0.72261613903606592
$ws.Range("BJ16").Value = 0.9578906035757061
$ws.Range("AV17").Value = 0.97009570216586383
$ws.Range("T18").Value = 0.85409516032612076
$ws.Range("R19").Value = 0.89672658565848085
$ws.Range("AN19").Value = 0.96741591517535097
$ws.Range("BC19").Value = 0.94551569314556805
$ws.Range("S20").Value = 0.99287210445768448
$ws.Range("V20").Value = 0.91213886375376685
$ws.Range("AC20").Value = 0.75207345912277201
$ws.Range("BK20").Value = 0.74142058160877611
$ws.Range("O21").Value = 0.84044430975469719
$ws.Range("AB21").Value = 0.97193850274084381
$ws.Range("AR21").Value = 0.87283109383564639
$ws.Range("Y22").Value = 0.77514875445542819
$ws.Range("AV22").Value = 0.80864349320828155
$ws.Range("O23").Value = 0.88641068953430524
$ws.Range("R23").Value = 0.82159452181616333
$ws.Range("AG23").Value = 0.98927994514252404
$ws.Range("BC23").Value = 0.7353682398473087
$ws.Range("W24").Value = 0.83661554293041607
$ws.Range("AQ24").Value = 0.65119899852414775
$ws.Range("BL24").Value = 0.87838985056890184
$ws.Range("BA25").Value = 0.78254095507673083
$ws.Range("AB26").Value = 0.92826948824501754
$ws.Range("AC27").Value = 0.98985795421185951
$ws.Range("AA28").Value = 0.9426621410009155
$ws.Range("AP28").Value = 0.62390801284578123
$ws.Range("AH29").Value = 0.96554567572919014
$ws.Range("AK29").Value = 0.94059455697466121
$ws.Range("AK31").Value = 0.92306193988050789
$ws.Range("AV31").Value = 0.85343043754124626
$ws.Range("BB31").Value = 0.70082889627073786
$ws.Range("A32").Value = 0.80301931061786946
$ws.Range("BK32").Value = 0.96261858944732748
$ws.Range("B33").Value = 0.92254006576397118
$ws.Range("AD34").Value = 0.91062598649110016
$ws.Range("AJ34").Value = 0.84046922361874044
$ws.Range("BK35").Value = 0.94523997010841643
$ws.Range("Y36").Value = 0.91627332531431427
$ws.Range("BK36").Value = 0.70958117162755929
$ws.Range("AI37").Value = 0.98250125232807162
$ws.Range("BN37").Value = 0.80537223510384059
$ws.Range("AM38").Value = 0.60906569935920629
$ws.Range("U39").Value = 0.96514072177198051
$ws.Range("AK39").Value = 0.77856605225775644
$ws.Range("AN39").Value = 0.98817780409563472
$ws.Range("AT39").Value = 0.89308908639272688
$ws.Range("K41").Value = 0.86101212977210873
$ws.Range("AF41").Value = 0.86730658681722239
$ws.Range("BN41").Value = 0.78598040524955926
$ws.Range("E42").Value = 0.85675640030282041
$ws.Range("H42").Value = 0.86706256649275826
$ws.Range("E43").Value = 0.94884253355759607
$ws.Range("K44").Value = 0.86677115201483756
$ws.Range("AB44").Value = 0.74713085718415329
$ws.Range("AD44").Value = 0.68572042406623446
$ws.Range("AH44").Value = 0.86090870385139895
$ws.Range("BI44").Value = 0.75766817162331279
$ws.Range("I45").Value = 0.95918956593122906
$ws.Range("AB45").Value = 0.59181247505603518
$ws.Range("AL45").Value = 0.84360554934918164
$ws.Range("BO45").Value = 0.99832564600732976
$ws.Range("S46").Value = 0.86569761977774073
$ws.Range("V46").Value = 0.84571978618530297
$ws.Range("BM47").Value = 0.92792312484257478
$ws.Range("P48").Value = 0.73324890702731338
$ws.Range("AJ49").Value = 0.88401793400965112
$ws.Range("AU49").Value = 0.80476072514140351
$ws.Range("BE49").Value = 0.82379034418406205
$ws.Range("A50").Value = 0.82169789371593582
$ws.Range("Z50").Value = 0.73574079293336925
$ws.Range("BJ50").Value = 0.88882341507714457
$ws.Range("A51").Value = 0.97885466688809863
$ws.Range("J51").Value = 0.87977116645278519
$ws.Range("AF51").Value = 0.83493934134768955
$ws.Range("AK51").Value = 0.9294351691608016
$ws.Range("AY52").Value = 0.79502272361554849
$ws.Range("AO53").Value = 0.98355415297505355
$ws.Range("AT53").Value = 0.84009202681403949
$ws.Range("AG54").Value = 0.90522954804649713
$ws.Range("BC54").Value = 0.78124038049875655
$ws.Range("H55").Value = 0.65253288402310905
$ws.Range("X55").Value = 0.86919225506915421
$ws.Range("AH55").Value = 0.95491857421502235
$ws.Range("AM56").Value = 0.73986401806688429
$ws.Range("AY56").Value = 0.88750973879325157
$ws.Range("E57").Value = 0.6538216008651162
$ws.Range("F58").Value = 0.9615317702004722
$ws.Range("V58").Value = 0.94967184082066214
$ws.Range("AZ58").Value = 0.71369124337585632
$ws.Range("AZ59").Value = 0.98487566340195243
$ws.Range("BE59").Value = 0.85367405791061235
$ws.Range("AI61").Value = 0.92645941357219852
$ws.Range("AA62").Value = 0.95663424672852027
$ws.Range("AT62").Value = 0.6698361146453008
$ws.Range("BH62").Value = 0.85981616376128223
$ws.Range("BL63").Value = 0.84811101263791211
$ws.Range("AS64").Value = 0.9666462733597746
$ws.Range("X65").Value = 0.94333611671655426
$ws.Range("AX65").Value = 0.80114627848212294
$ws.Range("BH65").Value = 0.99939418516841005
$ws.Range("BG66").Value = 0.96520113658919837
$ws.Range("M67").Value = 0.70615433245339643
$ws.Range("O67").Value = 0.85259468441514374
$ws.Range("BF67").Value = 0.99344752308411932
$ws.Range("D68").Value = 0.97039622529114755
$ws.Range("AO68").Value = 0.76904217382832996
